# Adds new meta type tkOverloadedFunction (and tkMetaRepo) rows to the
# "meta types" worksheet, matching the upstream commit
# "Added new meta type tkOverloadedFunction".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta types")

# --- Insert row for tkOverloadedFunction right after tkConstructor (row 25) ---
$ws.Rows.Item(26).Insert()
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 1).Value = "tkOverloadedFunction"
$ws.Cells.Item(26, 2).Value = 37
$ws.Cells.Item(26, 4).Value = "MetaCallable"
$ws.Cells.Item(26, 6).ClearContents()

# The enum values (column B) for tkDefaultArgsFunction .. tkMetaType (now on
# rows 27..33, since everything shifted down by one row) are contiguous
# integers that all bump by +1 because tkOverloadedFunction took the value
# previously held by tkDefaultArgsFunction.
for ($r = 27; $r -le 33; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 2).Value2 + 1
}

# --- Insert row for tkMetaRepo right after tkMetaType (old row 32, now row 33) ---
$ws.Rows.Item(34).Insert()
$ws.Cells.Item(34, 1).Value = "tkMetaRepo"
$ws.Cells.Item(34, 2).Value = 45
$ws.Cells.Item(34, 3).Value = "metapp::MetaRepo"

$ws.Range("D26").Select()
